# Re-run of the "write web table into Excel" data-provider tests in parallel:
#   - the three single-column sheets (Company / Contact / Country) and the two
#     whole-table sheets (nested-for / single-for) each get a second pass of
#     data written a few columns to the right of the first pass
#   - every "FineshedAt" timestamp cell (old AND new) gets refreshed, because
#     re-running the suite re-stamps completion time everywhere it writes.

$wb = $excel.ActiveWorkbook

$HEADER_FILL = 9868950   # same gray used by all existing column headers
$WIDTH_OFFSET = 0.8333333333333334   # this engine's char-width -> OOXML-width padding

function Write-MirroredTable {
    param(
        [string]$SheetName,
        [int]$SrcStartCol,
        [int]$ColCount,
        [string[]]$Headers,
        [string[][]]$Rows,
        [int]$FineshedAtCol,
        [string]$NewFineshedAtValue,
        [double[]]$ColWidths
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $dstStartCol = $SrcStartCol + $ColCount

    for ($c = 0; $c -lt $ColCount; $c++) {
        $dstCol = $dstStartCol + $c
        $cell1 = $ws.Cells.Item(1, $dstCol)
        $cell1.Value = $Headers[$c]
        $cell1.Interior.Color = $HEADER_FILL
    }

    for ($r = 0; $r -lt $Rows.Length; $r++) {
        $rowNum = $r + 2
        for ($c = 0; $c -lt $ColCount; $c++) {
            $val = $Rows[$r][$c]
            if ($val -ne $null) {
                $dstCol = $dstStartCol + $c
                $ws.Cells.Item($rowNum, $dstCol).Value = $val
            }
        }
    }

    $finCell = $ws.Cells.Item(4, $dstStartCol + $FineshedAtCol)
    $finCell.Value = $NewFineshedAtValue
    $ws.Rows.Item(4).AutoFit()

    for ($c = 0; $c -lt $ColCount; $c++) {
        $dstCol = $dstStartCol + $c
        $ws.Columns.Item($dstCol).ColumnWidth = ($ColWidths[$c] - $WIDTH_OFFSET)
    }
}

# ---------------------------------------------------------------------------
# Sheet: writeCompanyColumnIntoXcel  (A:B -> mirror into C:D)
# ---------------------------------------------------------------------------
$wsCompany = $wb.Worksheets.Item("writeCompanyColumnIntoXcel")
$wsCompany.Range("B4").Value = "1571855800605`nWed Oct 23 11:36:40 PDT 2019"
$wsCompany.Rows.Item(4).AutoFit()

Write-MirroredTable "writeCompanyColumnIntoXcel" 1 2 @("Company", "FineshedAt") @(
        @("Alfreds Futterkiste", $null),
        @("Centro comercial Moctezuma", $null),
        @("Ernst Handel", $null),
        @("Island Trading", $null),
        @("Laughing Bacchus Winecellars", $null),
        @("Magazzini Alimentari Riuniti", $null)
    ) 1 "1571855967414`nWed Oct 23 11:39:27 PDT 2019" @(28.36328125, 11.1328125)

# ---------------------------------------------------------------------------
# Sheet: writeContactColumnIntoXcel  (A:B -> mirror into C:D)
# ---------------------------------------------------------------------------
$wsContact = $wb.Worksheets.Item("writeContactColumnIntoXcel")
$wsContact.Range("B4").Value = "1571855802576`nWed Oct 23 11:36:42 PDT 2019"
$wsContact.Rows.Item(4).AutoFit()

Write-MirroredTable "writeContactColumnIntoXcel" 1 2 @("Contact", "FineshedAt") @(
        @("Maria Anders", $null),
        @("Francisco Chang", $null),
        @("Roland Mendel", $null),
        @("Helen Bennett", $null),
        @("Yoshi Tannamuri", $null),
        @("Giovanni Rovelli", $null)
    ) 1 "1571855969248`nWed Oct 23 11:39:29 PDT 2019" @(16.1328125, 11.1328125)

# ---------------------------------------------------------------------------
# Sheet: writeCountryColumnIntoXcel  (A:B -> mirror into C:D)
# ---------------------------------------------------------------------------
$wsCountry = $wb.Worksheets.Item("writeCountryColumnIntoXcel")
$wsCountry.Range("B4").Value = "1571855803803`nWed Oct 23 11:36:43 PDT 2019"
$wsCountry.Rows.Item(4).AutoFit()

Write-MirroredTable "writeCountryColumnIntoXcel" 1 2 @("Country", "FineshedAt") @(
        @("Germany", $null),
        @("Mexico", $null),
        @("Austria", $null),
        @("UK", $null),
        @("Canada", $null),
        @("Italy", $null)
    ) 1 "1571855971281`nWed Oct 23 11:39:31 PDT 2019" @(9.23828125, 11.1328125)

# ---------------------------------------------------------------------------
# Sheet: writeWholeTableNestedFor  (A:D -> mirror into E:H)
# ---------------------------------------------------------------------------
$wsNested = $wb.Worksheets.Item("writeWholeTableNestedFor")
$wsNested.Range("D4").Value = "1571855810696`nWed Oct 23 11:36:50 PDT 2019"
$wsNested.Rows.Item(4).AutoFit()

Write-MirroredTable "writeWholeTableNestedFor" 1 4 @("Company", "Contact", "Country", "FineshedAt") @(
        @("Alfreds Futterkiste", "Maria Anders", "Germany", $null),
        @("Centro comercial Moctezuma", "Francisco Chang", "Mexico", $null),
        @("Ernst Handel", "Roland Mendel", "Austria", $null),
        @("Island Trading", "Helen Bennett", "UK", $null),
        @("Laughing Bacchus Winecellars", "Yoshi Tannamuri", "Canada", $null),
        @("Magazzini Alimentari Riuniti", "Giovanni Rovelli", "Italy", $null)
    ) 3 "1571855997720`nWed Oct 23 11:39:57 PDT 2019" @(28.36328125, 16.1328125, 9.23828125, 11.1328125)

# ---------------------------------------------------------------------------
# Sheet: writeWholeTableSingleFor  (A:D -> mirror into E:H)
# ---------------------------------------------------------------------------
$wsSingle = $wb.Worksheets.Item("writeWholeTableSingleFor")
$wsSingle.Range("D4").Value = "1571855812034`nWed Oct 23 11:36:52 PDT 2019"
$wsSingle.Rows.Item(4).AutoFit()

Write-MirroredTable "writeWholeTableSingleFor" 1 4 @("Company", "Contact", "Country", "FineshedAt") @(
        @("Alfreds Futterkiste", "Maria Anders", "Germany", $null),
        @("Centro comercial Moctezuma", "Francisco Chang", "Mexico", $null),
        @("Ernst Handel", "Roland Mendel", "Austria", $null),
        @("Island Trading", "Helen Bennett", "UK", $null),
        @("Laughing Bacchus Winecellars", "Yoshi Tannamuri", "Canada", $null),
        @("Magazzini Alimentari Riuniti", "Giovanni Rovelli", "Italy", $null)
    ) 3 "1571856027433`nWed Oct 23 11:40:27 PDT 2019" @(28.36328125, 16.1328125, 9.23828125, 11.1328125)
